$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 562 (shifts old rows 562-569 down to 565-572),
# carrying over the existing formatting (e.g. the date style on column D).
$ws.Rows.Item(562).Resize(3).EntireRow.Insert()

# Populate the 3 newly inserted rows (562-564) with this week's price
# observations for "Vega Monumental Concepción" / Pimiento.
$newRows = @(
    @{ Row = 562; D = 45121; H = "Morrón rojo";   J = 100; K = 15000; L = 16000; M = 15500; N = "`$/caja 18 kilos"; O = "Provincia de Limarí";          P = 861;  Q = 18 },
    @{ Row = 563; D = 45121; H = "Zafiro rojo";    J = 100; K = 15000; L = 16000; M = 15500; N = "`$/caja 15 kilos"; O = "Región de Arica y Parinacota"; P = 1033; Q = 15 },
    @{ Row = 564; D = 45121; H = "Zafiro verde";   J = 100; K = 14000; L = 15000; M = 14500; N = "`$/caja 15 kilos"; O = "Región de Arica y Parinacota"; P = 967;  Q = 15 }
)

foreach ($rec in $newRows) {
    $r = $rec.Row
    $ws.Cells.Item($r, 1).Value = 11
    $ws.Cells.Item($r, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($r, 3).Value = "Bíobío"
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 5).Value = 8
    $ws.Cells.Item($r, 6).Value = 100112002
    $ws.Cells.Item($r, 7).Value = "Pimiento"
    $ws.Cells.Item($r, 8).Value = $rec.H
    $ws.Cells.Item($r, 9).Value = "Primera"
    $ws.Cells.Item($r, 10).Value = $rec.J
    $ws.Cells.Item($r, 11).Value = $rec.K
    $ws.Cells.Item($r, 12).Value = $rec.L
    $ws.Cells.Item($r, 13).Value = $rec.M
    $ws.Cells.Item($r, 14).Value = $rec.N
    $ws.Cells.Item($r, 15).Value = $rec.O
    $ws.Cells.Item($r, 16).Value = $rec.P
    $ws.Cells.Item($r, 17).Value = $rec.Q
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
